# Data Import Template.xlsx
# Insert a new "training_end" column between "training_date" (I) and
# "issue_date" (old J, now K), shifting issue_date/expiry_date/created_by
# one column to the right (J,K,L -> K,L,M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Copy the "training_date" column (I) and insert it as a new column at J so
# the new header inherits the same cell style (date-column text format) as
# its neighbours, then shift the old J:L columns (issue_date, expiry_date,
# created_by) one place to the right.
$ws.Columns.Item(9).Copy()
$ws.Columns.Item(10).Insert(-4161)          # -4161 = xlShiftToRight
$ws.Application.CutCopyMode = $false

# Relabel the newly inserted column.
$ws.Cells.Item(1, 10).Value = "training_end"

# Match the original column's width as closely as this host allows.
$ws.Columns.Item(10).ColumnWidth = 12.02

# Leftover UI state recorded in the saved file.
$ws.Range("L4").Select()
try {
    $excel.ActiveWindow.Width = 29040
    $excel.ActiveWindow.Height = 15720
} catch {}
